# Restructured and Added more documentation
#
# Adds a new bullet (sub-level) paragraph right after the paragraph that
# ends in "...which require the authorised login." and before the
# trailing "_GoBack" bookmark, i.e. the bookmark now belongs to the new
# paragraph instead of the old one. The new paragraph reuses the existing
# "List Paragraph" style / numbering definition (numId 1) but one level
# deeper (ilvl 1 instead of 0), and keeps the bold run formatting already
# used by the rest of that list.

$d = $word.ActiveDocument

# Locate the paragraph to split by searching for the tail of its text -
# this is more robust than hard-coding a paragraph index.
$find = $d.Content
$found = $find.Find.Execute(
    "which require the authorised login.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph text"
}

# Remember which paragraph the hit lives in, then re-fetch it from
# $d.Paragraphs (fresh/clean paragraph object) rather than from the Find
# range itself.
$anchorIndex = $find.Paragraphs.Item(1).Index
$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Split the paragraph right before its trailing paragraph mark (and
# before the bookmark that sits there), turning the bookmark + nothing
# else into a brand new, following paragraph.
$splitPos = $anchorPara.Range.End - 1
$splitPoint = $d.Range($splitPos, $splitPos)
$splitPoint.InsertBefore("`r")

# The newly created paragraph is the one right after the anchor.
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.InsertBefore("If you want to know how to login, please look at Android Application User Guide.")

# Demote it one level (ilvl 0 -> 1) within the same numbered list
# (numId 1); pStyle "List Paragraph" and bold run formatting are already
# inherited from the split, matching the rest of the list.
$newPara.Range.ListFormat.ListLevelNumber = 2
